$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.868.07"
$ws.Range("E2").Value = "  -3.98%  "

$ws.Range("D3").Value = "1.727.43"
$ws.Range("E3").Value = "  -2.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4922"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3512"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07243"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("E13").Value = "  -3.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.875"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").Value = "1.718.54"
$ws.Range("E15").Value = "  -2.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.809"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.18%  "

$ws.Range("E18").Value = "  -2.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06402"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("D23").Value = "26.932.84"
$ws.Range("E23").Value = "  -3.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.055"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.99%  "

$ws.Range("D28").Value = "1.914.47"
$ws.Range("E28").Value = "  -2.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.066"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.045"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09337"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("E33").Value = "  -2.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.364"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05890"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.424"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.44%  "

$ws.Range("E39").Value = "  -3.41%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1982"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.29%  "

$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9990"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5971"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.06%  "

$ws.Range("E43").Value = "  -6.48%  "

$ws.Range("E44").Value = "  -4.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.574"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5611"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.836"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06648"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.094"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.88%  "
